# Auto-generated edit script: update cached market-price / profit columns (H:N)
# across all 8 leve-profit worksheets, per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 933
$ws.Range("I6").Value = 399.5
$ws.Range("K6").Value = 1198.5
$ws.Range("M6").Value = -1086.5
$ws.Range("H9").Value = 222.8
$ws.Range("I9").Value = 203.5
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 203.5
$ws.Range("L9").Value = 300
$ws.Range("M9").Value = -34.5
$ws.Range("N9").Value = -638
$ws.Range("H12").Value = 803.1579
$ws.Range("I12").Value = 75.64286
$ws.Range("J12").Value = 2840.2
$ws.Range("K12").Value = 75.64286
$ws.Range("L12").Value = 2840.2
$ws.Range("M12").Value = 94.35714
$ws.Range("N12").Value = -3180.2
$ws.Range("H38").Value = 335.6
$ws.Range("I38").Value = 126
$ws.Range("J38").Value = 2222
$ws.Range("K38").Value = 378
$ws.Range("L38").Value = 6666
$ws.Range("M38").Value = -6
$ws.Range("N38").Value = -7410
$ws.Range("H53").Value = 1608.2667
$ws.Range("I53").Value = 307.16666
$ws.Range("J53").Value = 2475.6667
$ws.Range("K53").Value = 307.16666
$ws.Range("L53").Value = 2475.6667
$ws.Range("M53").Value = 329.83334
$ws.Range("N53").Value = -3749.6667
$ws.Range("H58").Value = 50
$ws.Range("I58").Value = 50
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 150
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 0
$ws.Range("N58").Value = $null
$ws.Range("H80").Value = 2332.4375
$ws.Range("I80").Value = 666
$ws.Range("K80").Value = 1998
$ws.Range("M80").Value = -1000
$ws.Range("H83").Value = 2332.4375
$ws.Range("I83").Value = 666
$ws.Range("K83").Value = 5994
$ws.Range("M83").Value = -1002
$ws.Range("H106").Value = 12919.533
$ws.Range("I106").Value = 2457
$ws.Range("J106").Value = 18150.8
$ws.Range("K106").Value = 2457
$ws.Range("L106").Value = 18150.8
$ws.Range("M106").Value = -1826
$ws.Range("N106").Value = -19412.8
$ws.Range("H121").Value = 2246.75
$ws.Range("J121").Value = 2246.75
$ws.Range("L121").Value = 6740.25
$ws.Range("N121").Value = -10234.25
$ws.Range("H137").Value = 3062.56
$ws.Range("I137").Value = 2054.7778
$ws.Range("K137").Value = 6164.3334
$ws.Range("M137").Value = -3614.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4581.28
$ws.Range("I61").Value = 2558.762
$ws.Range("K61").Value = 2558.762
$ws.Range("M61").Value = -2346.762
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96240
$ws.Range("H74").Value = 27780862
$ws.Range("I74").Value = 41670292
$ws.Range("K74").Value = 41670292
$ws.Range("M74").Value = -41669418
$ws.Range("H77").Value = 27780862
$ws.Range("I77").Value = 41670292
$ws.Range("K77").Value = 208351460
$ws.Range("M77").Value = -208347092
$ws.Range("H110").Value = 2135.1428
$ws.Range("I110").Value = 2135.1428
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 2135.1428
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -90.14280000000008
$ws.Range("N110").Value = $null
$ws.Range("H136").Value = 4581.28
$ws.Range("I136").Value = 2558.762
$ws.Range("K136").Value = 7676.286
$ws.Range("M136").Value = -5126.286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4487.273
$ws.Range("I86").Value = 2984.4443
$ws.Range("K86").Value = 2984.4443
$ws.Range("M86").Value = -1861.4443
$ws.Range("H89").Value = 4487.273
$ws.Range("I89").Value = 2984.4443
$ws.Range("K89").Value = 14922.2215
$ws.Range("M89").Value = -9306.2215
$ws.Range("H105").Value = 8939.611000000001
$ws.Range("I105").Value = 4600.3335
$ws.Range("K105").Value = 4600.3335
$ws.Range("M105").Value = -2853.3335
$ws.Range("H107").Value = 1451
$ws.Range("J107").Value = 999
$ws.Range("L107").Value = 999
$ws.Range("N107").Value = -4839
$ws.Range("H134").Value = 1577.4166
$ws.Range("I134").Value = 988.1429000000001
$ws.Range("K134").Value = 2964.4287
$ws.Range("M134").Value = -429.4287000000004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29388.65
$ws.Range("J31").Value = 69001.125
$ws.Range("L31").Value = 69001.125
$ws.Range("N31").Value = -69591.125
$ws.Range("H33").Value = 1550
$ws.Range("I33").Value = 1550
$ws.Range("K33").Value = 1550
$ws.Range("M33").Value = -1171
$ws.Range("H34").Value = 29388.65
$ws.Range("J34").Value = 69001.125
$ws.Range("L34").Value = 69001.125
$ws.Range("N34").Value = -69405.125
$ws.Range("H58").Value = 3947.8262
$ws.Range("I58").Value = 1742.3125
$ws.Range("J58").Value = 8989
$ws.Range("K58").Value = 1742.3125
$ws.Range("L58").Value = 8989
$ws.Range("M58").Value = -1539.3125
$ws.Range("N58").Value = -9395
$ws.Range("H74").Value = 138867.17
$ws.Range("J74").Value = 138867.17
$ws.Range("L74").Value = 138867.17
$ws.Range("N74").Value = -140615.17
$ws.Range("H77").Value = 138867.17
$ws.Range("J77").Value = 138867.17
$ws.Range("L77").Value = 416601.51
$ws.Range("N77").Value = -425337.51
$ws.Range("H136").Value = 3947.8262
$ws.Range("I136").Value = 1742.3125
$ws.Range("J136").Value = 8989
$ws.Range("K136").Value = 5226.9375
$ws.Range("L136").Value = 26967
$ws.Range("M136").Value = -2676.9375
$ws.Range("N136").Value = -32067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 402.66666
$ws.Range("I44").Value = 402.66666
$ws.Range("K44").Value = 1207.99998
$ws.Range("M44").Value = -809.9999800000001
$ws.Range("H137").Value = 3953.6924
$ws.Range("J137").Value = 4868.1113
$ws.Range("L137").Value = 14604.3339
$ws.Range("N137").Value = -24804.3339
$ws.Range("H140").Value = 3542.4285
$ws.Range("I140").Value = 3542.4285
$ws.Range("K140").Value = 10627.2855
$ws.Range("M140").Value = -5447.2855
$ws.Range("H141").Value = 4783.0586
$ws.Range("I141").Value = 2520
$ws.Range("K141").Value = 7560
$ws.Range("M141").Value = -2380

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14858.25
$ws.Range("J70").Value = 36271.57
$ws.Range("L70").Value = 36271.57
$ws.Range("N70").Value = -36811.57
$ws.Range("H73").Value = 14858.25
$ws.Range("J73").Value = 36271.57
$ws.Range("L73").Value = 36271.57
$ws.Range("N73").Value = -38143.57
$ws.Range("H80").Value = 4887.9565
$ws.Range("I80").Value = 2585.111
$ws.Range("J80").Value = 6368.357
$ws.Range("K80").Value = 2585.111
$ws.Range("L80").Value = 6368.357
$ws.Range("M80").Value = -1587.111
$ws.Range("N80").Value = -8364.357
$ws.Range("H83").Value = 4887.9565
$ws.Range("I83").Value = 2585.111
$ws.Range("J83").Value = 6368.357
$ws.Range("K83").Value = 12925.555
$ws.Range("L83").Value = 31841.785
$ws.Range("M83").Value = -7933.555
$ws.Range("N83").Value = -41825.785

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10625.25
$ws.Range("I22").Value = 3500
$ws.Range("K22").Value = 3500
$ws.Range("M22").Value = -3205
$ws.Range("H27").Value = 10625.25
$ws.Range("I27").Value = 3500
$ws.Range("K27").Value = 3500
$ws.Range("M27").Value = -3393
$ws.Range("H40").Value = 7939
$ws.Range("I40").Value = 7659.263
$ws.Range("K40").Value = 7659.263
$ws.Range("M40").Value = -7523.263
$ws.Range("H82").Value = 5688.864
$ws.Range("I82").Value = 3918.7144
$ws.Range("J82").Value = 8786.625
$ws.Range("K82").Value = 3918.7144
$ws.Range("L82").Value = 8786.625
$ws.Range("M82").Value = -3557.7144
$ws.Range("N82").Value = -9508.625
$ws.Range("H85").Value = 5688.864
$ws.Range("I85").Value = 3918.7144
$ws.Range("J85").Value = 8786.625
$ws.Range("K85").Value = 3918.7144
$ws.Range("L85").Value = 8786.625
$ws.Range("M85").Value = -2670.7144
$ws.Range("N85").Value = -11282.625
$ws.Range("H128").Value = 75000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 75000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 75000
$ws.Range("M128").Value = $null
$ws.Range("N128").Value = -84960
$ws.Range("H140").Value = 70836.25
$ws.Range("I140").Value = 75000
$ws.Range("J140").Value = 69448.336
$ws.Range("K140").Value = 75000
$ws.Range("L140").Value = 69448.336
$ws.Range("M140").Value = -69820
$ws.Range("N140").Value = -79808.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7476
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 7476
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
$ws.Range("H96").Value = 3964
$ws.Range("J96").Value = 5673.5
$ws.Range("L96").Value = 5673.5
$ws.Range("N96").Value = -8419.5
$ws.Range("H132").Value = 4097.879
$ws.Range("I132").Value = 3449.24
$ws.Range("K132").Value = 10347.72
$ws.Range("M132").Value = -7817.719999999999
$ws.Range("H136").Value = 1895.3103
$ws.Range("I136").Value = 1284.25
$ws.Range("K136").Value = 3852.75
$ws.Range("M136").Value = -1302.75
